$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns I1 (I0) and J1 (IF), styled like the other headers (H1)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-9: I column is always 1, J column mirrors H column
for ($r = 2; $r -le 9; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
